$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.711.60'
$ws.Range("E2").Value = '  +3.10%  '

$ws.Range("D3").Value = '1.861.67'
$ws.Range("E3").Value = '  +2.94%  '

$ws.Range("D4").Value = '1.036'
$ws.Range("E4").Value = '  +2.97%  '

$ws.Range("D5").Value = '323.99'
$ws.Range("E5").Value = '  +4.09%  '

$ws.Range("D6").Value = '1.033'
$ws.Range("E6").Value = '  +2.78%  '

$ws.Range("D7").Value = '0.4411'
$ws.Range("E7").Value = '  +3.02%  '

$ws.Range("D8").Value = '0.3807'
$ws.Range("E8").Value = '  +3.29%  '

$ws.Range("D9").Value = '0.07449'
$ws.Range("E9").Value = '  +3.27%  '

$ws.Range("D10").Value = '0.8824'
$ws.Range("E10").Value = '  +2.41%  '

$ws.Range("D11").Value = '21.72'
$ws.Range("E11").Value = '  +2.61%  '

$ws.Range("D12").Value = '1.871.32'
$ws.Range("E12").Value = '  -7.68%  '

$ws.Range("D13").Value = '5.551'
$ws.Range("E13").Value = '  +3.31%  '

$ws.Range("D14").Value = '6.739'
$ws.Range("E14").Value = '  +1.80%  '

$ws.Range("D15").Value = '0.07212'
$ws.Range("E15").Value = '  +4.67%  '

$ws.Range("D16").Value = '83.62'
$ws.Range("E16").Value = '  +3.79%  '

$ws.Range("D17").Value = '1.038'
$ws.Range("E17").Value = '  +3.22%  '

$ws.Range("D18").Value = '0.000009094'
$ws.Range("E18").Value = '  +2.12%  '

$ws.Range("D19").Value = '1.033'
$ws.Range("E19").Value = '  +2.88%  '

$ws.Range("D20").Value = '15.52'
$ws.Range("E20").Value = '  +2.38%  '

$ws.Range("D21").Value = '27.715.76'
$ws.Range("E21").Value = '  +2.87%  '

$ws.Range("D22").Value = '5.290'
$ws.Range("E22").Value = '  +2.09%  '

$ws.Range("E23").Value = '  +4.61%  '

$ws.Range("D24").Value = '158.57'
$ws.Range("E24").Value = '  +3.20%  '

$ws.Range("D25").Value = '1.930'
$ws.Range("E25").Value = '  +2.56%  '

$ws.Range("D26").Value = '18.82'
$ws.Range("E26").Value = '  +3.16%  '

$ws.Range("D27").Value = '1.996'
$ws.Range("E27").Value = '  +5.13%  '

$ws.Range("D28").Value = '5.312'
$ws.Range("E28").Value = '  +1.87%  '

$ws.Range("D29").Value = '117.49'
$ws.Range("E29").Value = '  +2.25%  '

$ws.Range("D30").Value = '0.09086'
$ws.Range("E30").Value = '  +1.79%  '

$ws.Range("D31").Value = '1.211'
$ws.Range("E31").Value = '  +4.72%  '

$ws.Range("D32").Value = '0.7661'
$ws.Range("E32").Value = '  +3.29%  '

$ws.Range("D33").Value = '4.576'
$ws.Range("E33").Value = '  +3.66%  '

$ws.Range("D34").Value = '2.897'
$ws.Range("E34").Value = '  +3.43%  '

$ws.Range("D35").Value = '1.034'
$ws.Range("E35").Value = '  +2.86%  '

$ws.Range("D36").Value = '1.159'
$ws.Range("E36").Value = '  +3.29%  '

$ws.Range("D37").Value = '0.01983'
$ws.Range("E37").Value = '  +3.34%  '

$ws.Range("D38").Value = '0.05347'
$ws.Range("E38").Value = '  +2.66%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.837'
$ws.Range("E39").Value = '  +3.73%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.5183'
$ws.Range("E40").Value = '  +2.17%  '

$ws.Range("D41").Value = '0.1688'
$ws.Range("E41").Value = '  +2.66%  '

$ws.Range("D42").Value = '6.813'
$ws.Range("E42").Value = '  +6.23%  '

$ws.Range("D43").Value = '8.641'
$ws.Range("E43").Value = '  +4.72%  '

$ws.Range("D44").Value = '109.48'
$ws.Range("E44").Value = '  +2.53%  '

$ws.Range("D45").Value = '10.59'
$ws.Range("E45").Value = '  +2.18%  '

$ws.Range("D46").Value = '1.722'
$ws.Range("E46").Value = '  +4.56%  '

$ws.Range("D47").Value = '0.4675'
$ws.Range("E47").Value = '  +2.38%  '

$ws.Range("D48").Value = '0.06426'
$ws.Range("E48").Value = '  +2.42%  '

$ws.Range("D49").Value = '1.864'
$ws.Range("E49").Value = '  +4.01%  '

$ws.Range("D50").Value = '39.72'
$ws.Range("E50").Value = '  +5.27%  '

$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").Value = '0.9337'
$ws.Range("E51").Value = '  +2.80%  '
